# Selenium Framework with TestNG & Logs generation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Run_Flag for TC_02 from "N" to "Y"
$ws.Range("A3").Value = "Y"

# Remove the now-unused trailing blank rows (4 and 5)
$ws.Rows("4:5").Delete()

# Move the active selection to C3
$ws.Range("C3").Select()
